$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8 (which holds A8=100) so that a
# new data point (75, 2.54) can be placed there, pushing 100/150/250 down.
$ws.Rows.Item(8).Insert()

# Update/molecular-weight distance data (A2:B11)
$data = @(
    @(10, 7.27),
    @(15, 6.51),
    @(20, 5.29),
    @(25, 4.7699999999999996),
    @(37, 4.08),
    @(50, 3.33),
    @(75, 2.54),
    @(100, 1.89),
    @(150, 0.95),
    @(250, 0.4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Update the active selection to match the saved view state
$ws.Range("D8").Select()
